# "ran model jan 26"
# Fill in the "Beat Vegas?" (column G) results for the games on 2021-01-25
# that were still pending, and append the newly-played/predicted games for
# 2021-01-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in column G ("Beat Vegas?") for rows 125-133 ---------------------
$beatVegas = @{
    125 = "No"
    126 = "No"
    127 = "No"
    128 = "No"
    129 = "No"
    130 = "Yes"
    131 = "No"
    132 = "No"
    133 = "Yes"
}

foreach ($row in $beatVegas.Keys) {
    $ws.Range("G$row").Value = $beatVegas[$row]
}

# --- Append new games for 2021-01-26 (rows 134-136) -------------------------
$newGames = @(
    @{ Row = 134; Date = 44222; Home = "ATL"; Away = "LAC"; Spread = -5;    Predicted = 11;     Diff = -16 }
    @{ Row = 135; Date = 44222; Home = "HOU"; Away = "WAS"; Spread = -3.5;  Predicted = 4.6;    Diff = -8.1 }
    @{ Row = 136; Date = 44222; Home = "UTA"; Away = "NYK"; Spread = -11;   Predicted = -22.2;  Diff = 11.2 }
)

foreach ($game in $newGames) {
    $r = $game.Row
    $ws.Range("A$r").Value = $game.Date
    $ws.Range("A$r").NumberFormat = "yyyy\-mm\-dd"
    $ws.Range("B$r").Value = $game.Home
    $ws.Range("C$r").Value = $game.Away
    $ws.Range("D$r").Value = $game.Spread
    $ws.Range("E$r").Value = $game.Predicted
    $ws.Range("F$r").Value = $game.Diff
}
